$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 40
$ws.Cells.Item($row, 1).Value = "Tommaso Bruschetti"
$ws.Cells.Item($row, 2).Value = "ELIA BATTISTI | U.S. Guarna"
$ws.Cells.Item($row, 3).Value = "Mattia Baldessarini | Shark Attack"
$ws.Cells.Item($row, 4).Value = "Sebastiano Zoller | CGB Gamberoni"
$ws.Cells.Item($row, 5).Value = "FEDERICO NICOLODI | U.S. Guarna"
$ws.Cells.Item($row, 6).Value = "Roberto Barozzi | Demobusters"
